$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 597.125
$ws.Range("I19").Value = 585.7143
$ws.Range("J19").Value = 606
$ws.Range("K19").Value = 585.7143
$ws.Range("L19").Value = 606
$ws.Range("M19").Value = -410.7143
$ws.Range("N19").Value = -956
$ws.Range("H33").Value = 251.71428
$ws.Range("I33").Value = 286
$ws.Range("J33").Value = 166
$ws.Range("K33").Value = 286
$ws.Range("L33").Value = 166
$ws.Range("M33").Value = -57
$ws.Range("N33").Value = -624
$ws.Range("H137").Value = 1785.5294
$ws.Range("I137").Value = 1330.2667
$ws.Range("J137").Value = 5200
$ws.Range("K137").Value = 3990.800099999999
$ws.Range("L137").Value = 15600
$ws.Range("M137").Value = -1440.800099999999
$ws.Range("N137").Value = -20700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1642.5416
$ws.Range("I2").Value = 1518.4375
$ws.Range("J2").Value = 1890.75
$ws.Range("K2").Value = 1518.4375
$ws.Range("L2").Value = 1890.75
$ws.Range("M2").Value = -1405.4375
$ws.Range("N2").Value = -2116.75
$ws.Range("H116").Value = 1642.5416
$ws.Range("I116").Value = 1518.4375
$ws.Range("J116").Value = 1890.75
$ws.Range("K116").Value = 1518.4375
$ws.Range("L116").Value = 1890.75
$ws.Range("M116").Value = 775.5625
$ws.Range("N116").Value = -6478.75
$ws.Range("H122").Value = 2528.9375
$ws.Range("J122").Value = 3012.5
$ws.Range("L122").Value = 9037.5
$ws.Range("N122").Value = -13937.5
$ws.Range("H132").Value = 3863.25
$ws.Range("I132").Value = 2966.6667
$ws.Range("J132").Value = 4401.2
$ws.Range("K132").Value = 8900.000100000001
$ws.Range("L132").Value = 13203.6
$ws.Range("M132").Value = -6370.000100000001
$ws.Range("N132").Value = -18263.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1642.5416
$ws.Range("I3").Value = 1518.4375
$ws.Range("J3").Value = 1890.75
$ws.Range("K3").Value = 1518.4375
$ws.Range("L3").Value = 1890.75
$ws.Range("M3").Value = -1404.4375
$ws.Range("N3").Value = -2118.75
$ws.Range("H20").Value = 144083.86
$ws.Range("I20").Value = 201117.4
$ws.Range("K20").Value = 201117.4
$ws.Range("M20").Value = -200870.4
$ws.Range("H94").Value = 1564
$ws.Range("I94").Value = 887.8
$ws.Range("K94").Value = 887.8
$ws.Range("M94").Value = -436.8
$ws.Range("H105").Value = 4467215.5
$ws.Range("I105").Value = 5955152.5
$ws.Range("K105").Value = 5955152.5
$ws.Range("M105").Value = -5953405.5
$ws.Range("H134").Value = 3673.077
$ws.Range("I134").Value = 3527.7778
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 10583.3334
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -8048.3334
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 39880
$ws.Range("J104").Value = 39880
$ws.Range("L104").Value = 39880
$ws.Range("N104").Value = -45122
$ws.Range("H105").Value = 1844.4615
$ws.Range("I105").Value = 1954.1428
$ws.Range("K105").Value = 1954.1428
$ws.Range("M105").Value = -207.1428000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4730.6665
$ws.Range("I5").Value = 7468
$ws.Range("J5").Value = 1993.3334
$ws.Range("K5").Value = 22404
$ws.Range("L5").Value = 5980.0002
$ws.Range("M5").Value = -22292
$ws.Range("N5").Value = -6204.0002
$ws.Range("H98").Value = 2093.3333
$ws.Range("I98").Value = 2308
$ws.Range("J98").Value = 1825
$ws.Range("K98").Value = 6924
$ws.Range("L98").Value = 5475
$ws.Range("M98").Value = -5426
$ws.Range("N98").Value = -8471
$ws.Range("H135").Value = 4730.6665
$ws.Range("I135").Value = 7468
$ws.Range("J135").Value = 1993.3334
$ws.Range("K135").Value = 67212
$ws.Range("L135").Value = 17940.0006
$ws.Range("M135").Value = -64677
$ws.Range("N135").Value = -23010.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5826.5654
$ws.Range("I70").Value = 5780.7334
$ws.Range("K70").Value = 5780.7334
$ws.Range("M70").Value = -5510.7334
$ws.Range("H73").Value = 5826.5654
$ws.Range("I73").Value = 5780.7334
$ws.Range("K73").Value = 5780.7334
$ws.Range("M73").Value = -4844.7334
$ws.Range("H80").Value = 3365.0588
$ws.Range("I80").Value = 2962.5
$ws.Range("J80").Value = 3722.889
$ws.Range("K80").Value = 2962.5
$ws.Range("L80").Value = 3722.889
$ws.Range("M80").Value = -1964.5
$ws.Range("N80").Value = -5718.889
$ws.Range("H83").Value = 3365.0588
$ws.Range("I83").Value = 2962.5
$ws.Range("J83").Value = 3722.889
$ws.Range("K83").Value = 14812.5
$ws.Range("L83").Value = 18614.445
$ws.Range("M83").Value = -9820.5
$ws.Range("N83").Value = -28598.445
$ws.Range("H122").Value = 4415.2
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4415.2
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 13245.6
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -18145.6
$ws.Range("H126").Value = 4390
$ws.Range("I126").Value = 4322.222
$ws.Range("K126").Value = 12966.666
$ws.Range("M126").Value = -10496.666
$ws.Range("H132").Value = 3097.7778
$ws.Range("I132").Value = 1633.3334
$ws.Range("J132").Value = 3830
$ws.Range("K132").Value = 4900.0002
$ws.Range("L132").Value = 11490
$ws.Range("M132").Value = -2370.0002
$ws.Range("N132").Value = -16550
$ws.Range("H134").Value = 198442
$ws.Range("J134").Value = 198442
$ws.Range("L134").Value = 595326
$ws.Range("N134").Value = -600396

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1059.091
$ws.Range("I16").Value = 950.8095
$ws.Range("K16").Value = 950.8095
$ws.Range("M16").Value = -780.8095
$ws.Range("H82").Value = 2392.0833
$ws.Range("I82").Value = 1257.4286
$ws.Range("J82").Value = 3980.6
$ws.Range("K82").Value = 1257.4286
$ws.Range("L82").Value = 3980.6
$ws.Range("M82").Value = -896.4286
$ws.Range("N82").Value = -4702.6
$ws.Range("H85").Value = 2392.0833
$ws.Range("I85").Value = 1257.4286
$ws.Range("J85").Value = 3980.6
$ws.Range("K85").Value = 1257.4286
$ws.Range("L85").Value = 3980.6
$ws.Range("M85").Value = -9.42859999999996
$ws.Range("N85").Value = -6476.6
$ws.Range("H122").Value = 64289428
$ws.Range("I122").Value = 41671000
$ws.Range("K122").Value = 125013000
$ws.Range("M122").Value = -125010550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2224
$ws.Range("H11").Value = 50000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = $null
$ws.Range("H13").Value = 3000
$ws.Range("I13").Value = 3666.6667
$ws.Range("K13").Value = 3666.6667
$ws.Range("M13").Value = -3526.6667
$ws.Range("H15").Value = 15000
$ws.Range("J15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("N15").Value = -15576
$ws.Range("H21").Value = 50000
$ws.Range("I21").Value = 50000
$ws.Range("K21").Value = 50000
$ws.Range("M21").Value = -49765
$ws.Range("H24").Value = 45008.332
$ws.Range("J24").Value = 44010
$ws.Range("L24").Value = 44010
$ws.Range("N24").Value = -44470
$ws.Range("H35").Value = 50000
$ws.Range("I35").Value = 50000
$ws.Range("K35").Value = 50000
$ws.Range("M35").Value = -49710
$ws.Range("H64").Value = 35000
$ws.Range("J64").Value = 35000
$ws.Range("L64").Value = 35000
$ws.Range("N64").Value = -35496
$ws.Range("H67").Value = 35000
$ws.Range("J67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("N67").Value = -36716
$ws.Range("H81").Value = 85665.92
$ws.Range("I81").Value = 123033.445
$ws.Range("J81").Value = 1589
$ws.Range("K81").Value = 246066.89
$ws.Range("L81").Value = 3178
$ws.Range("M81").Value = -245005.89
$ws.Range("N81").Value = -5300
$ws.Range("H84").Value = 85665.92
$ws.Range("I84").Value = 123033.445
$ws.Range("J84").Value = 1589
$ws.Range("K84").Value = 1230334.45
$ws.Range("L84").Value = 15890
$ws.Range("M84").Value = -1225030.45
$ws.Range("N84").Value = -26498
$ws.Range("H132").Value = 3766.9167
$ws.Range("I132").Value = 3260
$ws.Range("J132").Value = 4129
$ws.Range("K132").Value = 9780
$ws.Range("L132").Value = 12387
$ws.Range("M132").Value = -7250
$ws.Range("N132").Value = -17447
